$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" '29.094.02'
Set-TextValue $ws "E2" '  -0.04%  '
Set-TextValue $ws "D3" '1.821.78'
Set-TextValue $ws "E3" '  -0.62%  '
Set-TextValue $ws "D4" '0.9977'
Set-TextValue $ws "E4" '  -0.42%  '
Set-TextValue $ws "D5" '241.33'
Set-TextValue $ws "E5" '  -0.86%  '
Set-TextValue $ws "D6" '0.6146'
Set-TextValue $ws "E6" '  -2.10%  '
Set-TextValue $ws "D7" '0.9996'
Set-TextValue $ws "E7" '  -0.28%  '
Set-TextValue $ws "D8" '0.07321'
Set-TextValue $ws "E8" '  -2.18%  '
Set-TextValue $ws "D9" '0.2884'
Set-TextValue $ws "E9" '  -1.38%  '
Set-TextValue $ws "D10" '22.90'
Set-TextValue $ws "E10" '  -1.38%  '
Set-TextValue $ws "D11" '0.07657'
Set-TextValue $ws "E11" '  -0.36%  '
Set-TextValue $ws "D12" '1.816.30'
Set-TextValue $ws "E12" '  -0.95%  '
Set-TextValue $ws "D13" '4.944'
Set-TextValue $ws "E13" '  -1.21%  '
Set-TextValue $ws "D14" '0.6584'
Set-TextValue $ws "E14" '  -1.32%  '
Set-TextValue $ws "D15" '81.65'
Set-TextValue $ws "E15" '  -1.25%  '
Set-TextValue $ws "D16" '0.000008946'
Set-TextValue $ws "E16" '  -4.28%  '
Set-TextValue $ws "D17" '5.817'
Set-TextValue $ws "E17" '  -2.78%  '
Set-TextValue $ws "D18" '29.061.44'
Set-TextValue $ws "E18" '  -0.16%  '
Set-TextValue $ws "D19" '2.066.10'
Set-TextValue $ws "E19" '  -0.39%  '
Set-TextValue $ws "D20" '237.48'
Set-TextValue $ws "E20" '  +6.44%  '
Set-TextValue $ws "D21" '12.41'
Set-TextValue $ws "E21" '  -1.27%  '
Set-TextValue $ws "D22" '0.9992'
Set-TextValue $ws "E22" '  -0.43%  '
Set-TextValue $ws "D23" '7.103'
Set-TextValue $ws "E23" '  +0.11%  '
Set-TextValue $ws "D24" '0.9998'
Set-TextValue $ws "E24" '  -0.28%  '
Set-TextValue $ws "D25" '157.31'
Set-TextValue $ws "E25" '  -1.59%  '
Set-TextValue $ws "D26" '0.1405'
Set-TextValue $ws "E26" '  +1.19%  '
Set-TextValue $ws "D27" '8.414'
Set-TextValue $ws "E27" '  -0.80%  '
Set-TextValue $ws "D28" '17.59'
Set-TextValue $ws "E28" '  -1.67%  '
Set-TextValue $ws "D29" '1.481'
Set-TextValue $ws "E29" '  -1.29%  '
Set-TextValue $ws "D30" '0.05545'
Set-TextValue $ws "E30" '  -3.17%  '
Set-TextValue $ws "D31" '4.082'
Set-TextValue $ws "E31" '  +0.12%  '
Set-TextValue $ws "D32" '4.083'
Set-TextValue $ws "E32" '  -1.55%  '
Set-TextValue $ws "D33" '1.199'
Set-TextValue $ws "E33" '  -0.60%  '
Set-TextValue $ws "D34" '1.823'
Set-TextValue $ws "E34" '  -0.38%  '
Set-TextValue $ws "D35" '0.7338'
Set-TextValue $ws "E35" '  -1.05%  '
Set-TextValue $ws "D36" '1.129'
Set-TextValue $ws "E36" '  -0.90%  '
Set-TextValue $ws "D37" '2.608'
Set-TextValue $ws "E37" '  -2.37%  '
Set-TextValue $ws "D38" '2.829'
Set-TextValue $ws "E38" '  +2.30%  '
Set-TextValue $ws "D39" '1.206.98'
Set-TextValue $ws "E39" '  -0.49%  '
Set-TextValue $ws "D40" '0.01755'
Set-TextValue $ws "E40" '  -1.29%  '
Set-TextValue $ws "D41" '6.349'
Set-TextValue $ws "E41" '  -2.58%  '
Set-TextValue $ws "D42" '0.8949'
Set-TextValue $ws "E42" '  +0.58%  '
Set-TextValue $ws "D43" '1.000'
Set-TextValue $ws "E43" '  -0.17%  '
Set-TextValue $ws "D44" '101.04'
Set-TextValue $ws "E44" '  -0.84%  '
Set-TextValue $ws "D45" '1.970.58'
Set-TextValue $ws "E45" '  -0.49%  '
Set-TextValue $ws "D46" '64.53'
Set-TextValue $ws "E46" '  -1.46%  '
Set-TextValue $ws "D47" '0.5079'
Set-TextValue $ws "E47" '  -0.33%  '
Set-TextValue $ws "D48" '0.00000000117'
Set-TextValue $ws "E48" '  -6.00%  '
Set-TextValue $ws "D49" '0.3993'
Set-TextValue $ws "E49" '  -1.61%  '
Set-TextValue $ws "D50" '8.979'
Set-TextValue $ws "E50" '  -0.31%  '
Set-TextValue $ws "D51" '0.05745'
Set-TextValue $ws "E51" '  -1.37%  '
